$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the date/time formatting of row 24 down into the two new rows
# (25 and 26) before touching any values, so the new cells pick up the
# same cell styles already used by the existing data (A: date, B: time).
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("A26").PasteSpecial(-4122)

$ws.Range("B24").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("B26").PasteSpecial(-4122)

# Fix the date in A24 (was a typo: 41369 -> 41583)
$ws.Range("A24").Value = 41583

# New row 25
$ws.Range("A25").Value = 41584
$ws.Range("B25").Value = 0.09027777777777778

# New row 26
$ws.Range("A26").Value = 41588
$ws.Range("B26").Value = 0.09722222222222222

# Update selection to match the new active cell after data entry
$ws.Range("C26").Select()
